# Fulfillment requests report template:
# - Insert a new column before column C and move "Request Status" there
#   (it used to be the last column, Z).
# - Add a new trailing "Asset Status" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column B's width so the freshly inserted column C can match it
# (mirrors Excel's own "inherit formatting from the column to the left").
$colBWidth = $ws.Range("B1").EntireColumn.ColumnWidth

# Insert a blank column before column C; everything from C.. shifts to D..
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").EntireColumn.ColumnWidth = $colBWidth

# "Request Status" used to live in column Z; after the insert above it is
# now in column AA. Move its header text into the freshly inserted column C
# and replace the old slot with the new "Asset Status" header.
$ws.Range("C1").Value = $ws.Range("AA1").Value2
$ws.Range("AA1").Value = "Asset Status"

# Re-apply the autofilter so its range covers the new A1:AA1 header row.
$ws.AutoFilterMode = $false
$ws.Range("A1:AA1").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name (driven by the autofilter) still
# points at the old A1:Z1 range; repoint it at the new extent.
$wb.Names.Item(1).RefersTo = '=Data!$A$1:$AA$1'
